$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change sheet: update data rows 2-9 with new values ---
$ws1.Cells.Item(2, 1).Value = "HDFCBANK"
$ws1.Cells.Item(2, 2).Value = 944.7
$ws1.Cells.Item(2, 3).Value = 953
$ws1.Cells.Item(2, 4).Value = 941.4
$ws1.Cells.Item(2, 5).Value = 948
$ws1.Cells.Item(2, 6).Value = 21928849
$ws1.Cells.Item(2, 7).Value = 52884712
$ws1.Cells.Item(2, 8).Value = -0.5853461582621458
$ws1.Cells.Item(2, 9).Value = "HDFCBANK"
$ws1.Cells.Item(3, 1).Value = "MAXHEALTH"
$ws1.Cells.Item(3, 2).Value = 1030
$ws1.Cells.Item(3, 3).Value = 1035.1
$ws1.Cells.Item(3, 4).Value = 1019.1
$ws1.Cells.Item(3, 5).Value = 1030
$ws1.Cells.Item(3, 6).Value = 2117653
$ws1.Cells.Item(3, 7).Value = 4484760
$ws1.Cells.Item(3, 8).Value = -0.527811298709407
$ws1.Cells.Item(3, 9).Value = "MAXHEALTH"
$ws1.Cells.Item(4, 1).Value = "GODREJCP"
$ws1.Cells.Item(4, 2).Value = 1246.7
$ws1.Cells.Item(4, 3).Value = 1249.7
$ws1.Cells.Item(4, 4).Value = 1230.6
$ws1.Cells.Item(4, 5).Value = 1236.4
$ws1.Cells.Item(4, 6).Value = 680283
$ws1.Cells.Item(4, 7).Value = 1553698
$ws1.Cells.Item(4, 8).Value = -0.5621523616558688
$ws1.Cells.Item(4, 9).Value = "GODREJCP"
$ws1.Cells.Item(5, 1).Value = "NAUKRI"
$ws1.Cells.Item(5, 2).Value = 1366
$ws1.Cells.Item(5, 3).Value = 1368
$ws1.Cells.Item(5, 4).Value = 1333
$ws1.Cells.Item(5, 5).Value = 1340.3
$ws1.Cells.Item(5, 6).Value = 1070292
$ws1.Cells.Item(5, 7).Value = 2475123
$ws1.Cells.Item(5, 8).Value = -0.567580277828617
$ws1.Cells.Item(5, 9).Value = "NAUKRI"
$ws1.Cells.Item(6, 1).Value = "LTIM"
$ws1.Cells.Item(6, 2).Value = 6107
$ws1.Cells.Item(6, 3).Value = 6134
$ws1.Cells.Item(6, 4).Value = 6000.5
$ws1.Cells.Item(6, 5).Value = 6004
$ws1.Cells.Item(6, 6).Value = 146676
$ws1.Cells.Item(6, 7).Value = 292945
$ws1.Cells.Item(6, 8).Value = -0.4993053303521139
$ws1.Cells.Item(6, 9).Value = "LTIM"
$ws1.Cells.Item(7, 1).Value = "GMRAIRPORT"
$ws1.Cells.Item(7, 2).Value = 104.11
$ws1.Cells.Item(7, 3).Value = 104.99
$ws1.Cells.Item(7, 4).Value = 102.15
$ws1.Cells.Item(7, 5).Value = 102.48
$ws1.Cells.Item(7, 6).Value = 5256346
$ws1.Cells.Item(7, 7).Value = 11690565
$ws1.Cells.Item(7, 8).Value = -0.5503770775835043
$ws1.Cells.Item(7, 9).Value = "GMRAIRPORT"
$ws1.Cells.Item(8, 1).Value = "MANAPPURAM"
$ws1.Cells.Item(8, 2).Value = 319.9
$ws1.Cells.Item(8, 3).Value = 320.2
$ws1.Cells.Item(8, 4).Value = 306.4
$ws1.Cells.Item(8, 5).Value = 308.5
$ws1.Cells.Item(8, 6).Value = 4155251
$ws1.Cells.Item(8, 7).Value = 8374754
$ws1.Cells.Item(8, 8).Value = -0.5038360529754068
$ws1.Cells.Item(8, 9).Value = "MANAPPURAM"
$ws1.Cells.Item(9, 1).Value = "PPLPHARMA"
$ws1.Cells.Item(9, 2).Value = 180.81
$ws1.Cells.Item(9, 3).Value = 181.48
$ws1.Cells.Item(9, 4).Value = 172
$ws1.Cells.Item(9, 5).Value = 172.35
$ws1.Cells.Item(9, 6).Value = 3119826
$ws1.Cells.Item(9, 7).Value = 6506808
$ws1.Cells.Item(9, 8).Value = -0.5205289598217744
$ws1.Cells.Item(9, 9).Value = "PPLPHARMA"

# Remove now-unused rows 10-12 (sheet shrinks from 12 to 8 data rows)
$ws1.Range("A10:I12").EntireRow.Delete()

# --- Pos_Change sheet: update data rows 2-13 with new values ---
$ws2.Cells.Item(2, 1).Value = "BAJAJ-AUTO"
$ws2.Cells.Item(2, 2).Value = 9799.5
$ws2.Cells.Item(2, 3).Value = 9888
$ws2.Cells.Item(2, 4).Value = 9725
$ws2.Cells.Item(2, 5).Value = 9738
$ws2.Cells.Item(2, 6).Value = 534086
$ws2.Cells.Item(2, 7).Value = 334802
$ws2.Cells.Item(2, 8).Value = 0.5952294191790969
$ws2.Cells.Item(2, 9).Value = "BAJAJ-AUTO"
$ws2.Cells.Item(3, 1).Value = "HDFCLIFE"
$ws2.Cells.Item(3, 2).Value = 773
$ws2.Cells.Item(3, 3).Value = 774.35
$ws2.Cells.Item(3, 4).Value = 752.65
$ws2.Cells.Item(3, 5).Value = 759.9
$ws2.Cells.Item(3, 6).Value = 1993731
$ws2.Cells.Item(3, 7).Value = 1349921
$ws2.Cells.Item(3, 8).Value = 0.4769242051942299
$ws2.Cells.Item(3, 9).Value = "HDFCLIFE"
$ws2.Cells.Item(4, 1).Value = "GRASIM"
$ws2.Cells.Item(4, 2).Value = 2821
$ws2.Cells.Item(4, 3).Value = 2831.9
$ws2.Cells.Item(4, 4).Value = 2787.1
$ws2.Cells.Item(4, 5).Value = 2787.1
$ws2.Cells.Item(4, 6).Value = 437229
$ws2.Cells.Item(4, 7).Value = 287536
$ws2.Cells.Item(4, 8).Value = 0.5206061154081576
$ws2.Cells.Item(4, 9).Value = "GRASIM"
$ws2.Cells.Item(5, 1).Value = "RELIANCE"
$ws2.Cells.Item(5, 2).Value = 1500
$ws2.Cells.Item(5, 3).Value = 1503.9
$ws2.Cells.Item(5, 4).Value = 1468.8
$ws2.Cells.Item(5, 5).Value = 1470.7
$ws2.Cells.Item(5, 6).Value = 16518684
$ws2.Cells.Item(5, 7).Value = 11199340
$ws2.Cells.Item(5, 8).Value = 0.4749694178406942
$ws2.Cells.Item(5, 9).Value = "RELIANCE"
$ws2.Cells.Item(6, 1).Value = "ADANIENT"
$ws2.Cells.Item(6, 2).Value = 2275
$ws2.Cells.Item(6, 3).Value = 2275.9
$ws2.Cells.Item(6, 4).Value = 2205
$ws2.Cells.Item(6, 5).Value = 2209.3
$ws2.Cells.Item(6, 6).Value = 792800
$ws2.Cells.Item(6, 7).Value = 529451
$ws2.Cells.Item(6, 8).Value = 0.497400137123171
$ws2.Cells.Item(6, 9).Value = "ADANIENT"
$ws2.Cells.Item(7, 1).Value = "MOTHERSON"
$ws2.Cells.Item(7, 2).Value = 118
$ws2.Cells.Item(7, 3).Value = 118.94
$ws2.Cells.Item(7, 4).Value = 117.25
$ws2.Cells.Item(7, 5).Value = 117.27
$ws2.Cells.Item(7, 6).Value = 13950281
$ws2.Cells.Item(7, 7).Value = 9250970
$ws2.Cells.Item(7, 8).Value = 0.5079803523306204
$ws2.Cells.Item(7, 9).Value = "MOTHERSON"
$ws2.Cells.Item(8, 1).Value = "BRITANNIA"
$ws2.Cells.Item(8, 2).Value = 6155
$ws2.Cells.Item(8, 3).Value = 6158.5
$ws2.Cells.Item(8, 4).Value = 5945.5
$ws2.Cells.Item(8, 5).Value = 6022
$ws2.Cells.Item(8, 6).Value = 641813
$ws2.Cells.Item(8, 7).Value = 437361
$ws2.Cells.Item(8, 8).Value = 0.4674673782070189
$ws2.Cells.Item(8, 9).Value = "BRITANNIA"
$ws2.Cells.Item(9, 1).Value = "IDFCFIRSTB"
$ws2.Cells.Item(9, 2).Value = 84.39
$ws2.Cells.Item(9, 3).Value = 86.45
$ws2.Cells.Item(9, 4).Value = 83.38
$ws2.Cells.Item(9, 5).Value = 85.88
$ws2.Cells.Item(9, 6).Value = 62702609
$ws2.Cells.Item(9, 7).Value = 42923311
$ws2.Cells.Item(9, 8).Value = 0.4608055049620939
$ws2.Cells.Item(9, 9).Value = "IDFCFIRSTB"
$ws2.Cells.Item(10, 1).Value = "CUMMINSIND"
$ws2.Cells.Item(10, 2).Value = 4148
$ws2.Cells.Item(10, 3).Value = 4235
$ws2.Cells.Item(10, 4).Value = 4096
$ws2.Cells.Item(10, 5).Value = 4126
$ws2.Cells.Item(10, 6).Value = 1117650
$ws2.Cells.Item(10, 7).Value = 778593
$ws2.Cells.Item(10, 8).Value = 0.4354739896197372
$ws2.Cells.Item(10, 9).Value = "CUMMINSIND"
$ws2.Cells.Item(11, 1).Value = "ATGL"
$ws2.Cells.Item(11, 2).Value = 572
$ws2.Cells.Item(11, 3).Value = 577.8
$ws2.Cells.Item(11, 4).Value = 561.05
$ws2.Cells.Item(11, 5).Value = 561.7
$ws2.Cells.Item(11, 6).Value = 526321
$ws2.Cells.Item(11, 7).Value = 344248
$ws2.Cells.Item(11, 8).Value = 0.5289006762566522
$ws2.Cells.Item(11, 9).Value = "ATGL"
$ws2.Cells.Item(12, 1).Value = "MARICO"
$ws2.Cells.Item(12, 2).Value = 773.55
$ws2.Cells.Item(12, 3).Value = 773.55
$ws2.Cells.Item(12, 4).Value = 756.35
$ws2.Cells.Item(12, 5).Value = 760.5
$ws2.Cells.Item(12, 6).Value = 3203911
$ws2.Cells.Item(12, 7).Value = 2185434
$ws2.Cells.Item(12, 8).Value = 0.4660296307278097
$ws2.Cells.Item(12, 9).Value = "MARICO"
$ws2.Cells.Item(13, 1).Value = "UNIONBANK"
$ws2.Cells.Item(13, 2).Value = 166.28
$ws2.Cells.Item(13, 3).Value = 166.92
$ws2.Cells.Item(13, 4).Value = 160.25
$ws2.Cells.Item(13, 5).Value = 161.27
$ws2.Cells.Item(13, 6).Value = 13384492
$ws2.Cells.Item(13, 7).Value = 9352445
$ws2.Cells.Item(13, 8).Value = 0.4311222359500644
$ws2.Cells.Item(13, 9).Value = "UNIONBANK"
